# Weekly update: insert a new observation for "Brócoli" (Hortaliza,
# Macroferia Regional de Talca) as row 70, pushing the existing rows
# 70-160 down to 71-161 (so the previously-last row, 160, becomes 161).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 70; everything below (old rows 70..160)
# shifts down by one (to 71..161), carrying its formatting with it.
$ws.Rows("70:70").Insert()

# Populate the newly inserted row 70 with the new weekly record.
$ws.Range("A70").Value = 5
$ws.Range("B70").Value = "Macroferia Regional de Talca"
$ws.Range("C70").Value = "Maule"
$ws.Range("D70").Value = 44413
$ws.Range("E70").Value = 7
$ws.Range("F70").Value = 100112023
$ws.Range("G70").Value = "Brócoli"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Segunda"
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 400
$ws.Range("L70").Value = 400
$ws.Range("M70").Value = 400
$ws.Range("N70").Value = '$/unidad'
$ws.Range("O70").Value = "Región del Maule"
$ws.Range("P70").Value = 400
$ws.Range("Q70").Value = 1
$ws.Range("R70").Value = "Hortaliza"

# Match the date column's existing number format for the new row.
$ws.Range("D70").NumberFormat = $ws.Range("D71").NumberFormat
